$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Replace the "Good Morning" greeting for rule R10 with "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active cell selection left behind by the edit
$ws.Activate()
$ws.Range("E8").Select()
